$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 46; this shifts the existing rows 46-51 down to 47-52,
# carrying along their formatting (including the date style on column D).
$ws.Rows.Item(46).Insert()

# Populate the newly inserted row 46 with the new weekly price report.
# Columns A, B, C, E, F, G, H, I, N, O, Q, R are constant for this market/product
# across the whole block of rows, so write the same literal values used elsewhere.
$ws.Cells.Item(46, 1).Value = 1
$ws.Cells.Item(46, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(46, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(46, 4).Value = 44474
$ws.Cells.Item(46, 5).Value = 15
$ws.Cells.Item(46, 6).Value = 100112040
$ws.Cells.Item(46, 7).Value = "Cilantro"
$ws.Cells.Item(46, 8).Value = "Sin especificar"
$ws.Cells.Item(46, 9).Value = "Primera"
$ws.Cells.Item(46, 10).Value = 300
$ws.Cells.Item(46, 11).Value = 700
$ws.Cells.Item(46, 12).Value = 800
$ws.Cells.Item(46, 13).Value = 750
$ws.Cells.Item(46, 14).Value = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(46, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(46, 16).Value = 375
$ws.Cells.Item(46, 17).Value = 2
$ws.Cells.Item(46, 18).Value = "Hortaliza"
